# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on several leve rows across multiple crafting-job
# sheets with newly fetched Universalis price snapshots.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3890
$ws.Range("I64").Value = 3623.077
$ws.Range("J64").Value = 4385.7144
$ws.Range("K64").Value = 3623.077
$ws.Range("L64").Value = 4385.7144
$ws.Range("M64").Value = -3375.077
$ws.Range("N64").Value = -4881.7144

$ws.Range("H67").Value = 3890
$ws.Range("I67").Value = 3623.077
$ws.Range("J67").Value = 4385.7144
$ws.Range("K67").Value = 3623.077
$ws.Range("L67").Value = 4385.7144
$ws.Range("M67").Value = -2765.077
$ws.Range("N67").Value = -6101.7144

$ws.Range("H74").Value = 4190.909
$ws.Range("I74").Value = 3850
$ws.Range("J74").Value = 4600
$ws.Range("K74").Value = 3850
$ws.Range("L74").Value = 4600
$ws.Range("M74").Value = -2914
$ws.Range("N74").Value = -6472

$ws.Range("H77").Value = 4190.909
$ws.Range("I77").Value = 3850
$ws.Range("J77").Value = 4600
$ws.Range("K77").Value = 19250
$ws.Range("L77").Value = 23000
$ws.Range("M77").Value = -14570
$ws.Range("N77").Value = -32360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10118.624
$ws.Range("I32").Value = 7543.073
$ws.Range("K32").Value = 7543.073
$ws.Range("M32").Value = -7256.073

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5877.4443
$ws.Range("I82").Value = 2076.75
$ws.Range("J82").Value = 36283
$ws.Range("K82").Value = 2076.75
$ws.Range("L82").Value = 36283
$ws.Range("M82").Value = -1693.75
$ws.Range("N82").Value = -37049

$ws.Range("H85").Value = 5877.4443
$ws.Range("I85").Value = 2076.75
$ws.Range("J85").Value = 36283
$ws.Range("K85").Value = 2076.75
$ws.Range("L85").Value = 36283
$ws.Range("M85").Value = -750.75
$ws.Range("N85").Value = -38935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6011.8237
$ws.Range("I31").Value = 1110.7273
$ws.Range("J31").Value = 14997.167
$ws.Range("K31").Value = 1110.7273
$ws.Range("L31").Value = 14997.167
$ws.Range("M31").Value = -815.7273
$ws.Range("N31").Value = -15587.167

$ws.Range("H34").Value = 6011.8237
$ws.Range("I34").Value = 1110.7273
$ws.Range("J34").Value = 14997.167
$ws.Range("K34").Value = 1110.7273
$ws.Range("L34").Value = 14997.167
$ws.Range("M34").Value = -908.7273
$ws.Range("N34").Value = -15401.167

$ws.Range("H132").Value = 3013.2222
$ws.Range("J132").Value = 4249.5
$ws.Range("L132").Value = 12748.5
$ws.Range("N132").Value = -17808.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1507.6923
$ws.Range("I4").Value = 720
$ws.Range("K4").Value = 2160
$ws.Range("M4").Value = -2048

$ws.Range("H68").Value = 13801.75
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 15659.143
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 46977.429
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -48599.429

$ws.Range("H71").Value = 13801.75
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 15659.143
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 140932.287
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -149044.287

$ws.Range("H103").Value = 885.1429000000001
$ws.Range("I103").Value = 280
$ws.Range("J103").Value = 1339
$ws.Range("K103").Value = 840
$ws.Range("L103").Value = 4017
$ws.Range("M103").Value = 39
$ws.Range("N103").Value = -5775

$ws.Range("H112").Value = 4338.3335
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4338.3335
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 13015.0005
$ws.Range("N112").Value = -15231.0005
$ws.Range("M112").ClearContents()

$ws.Range("H114").Value = 780.46155
$ws.Range("J114").Value = 838.8570999999999
$ws.Range("L114").Value = 2516.5713
$ws.Range("N114").Value = -9024.5713

$ws.Range("H122").Value = 1140.3684
$ws.Range("J122").Value = 1521.7826
$ws.Range("L122").Value = 13696.0434
$ws.Range("N122").Value = -18596.0434

$ws.Range("H125").Value = 2635.1
$ws.Range("I125").Value = 887.75
$ws.Range("J125").Value = 3800
$ws.Range("K125").Value = 2663.25
$ws.Range("L125").Value = 11400
$ws.Range("M125").Value = 2256.75
$ws.Range("N125").Value = -21240

$ws.Range("H129").Value = 2423.9412
$ws.Range("J129").Value = 1722
$ws.Range("L129").Value = 5166
$ws.Range("N129").Value = -15166

$ws.Range("H131").Value = 12168.429
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 12168.429
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 36505.287
$ws.Range("N131").Value = -46585.287
$ws.Range("M131").ClearContents()

$ws.Range("H137").Value = 38465820
$ws.Range("I137").Value = 45457144
$ws.Range("K137").Value = 136371432
$ws.Range("M137").Value = -136366332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 40566.17
$ws.Range("I132").Value = 57468.168
$ws.Range("K132").Value = 172404.504
$ws.Range("M132").Value = -169874.504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1069.1538
$ws.Range("J46").Value = 1029.9
$ws.Range("L46").Value = 1029.9
$ws.Range("N46").Value = -1405.9

$ws.Range("H55").Value = 236242
$ws.Range("I55").Value = 400641.8
$ws.Range("J55").Value = 1385.1428
$ws.Range("K55").Value = 400641.8
$ws.Range("L55").Value = 1385.1428
$ws.Range("M55").Value = -400468.8
$ws.Range("N55").Value = -1731.1428

$ws.Range("H132").Value = 5252.909
$ws.Range("I132").Value = 5688.75
$ws.Range("J132").Value = 4090.6667
$ws.Range("K132").Value = 17066.25
$ws.Range("L132").Value = 12272.0001
$ws.Range("M132").Value = -14536.25
$ws.Range("N132").Value = -17332.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 18999
$ws.Range("J54").Value = 18999
$ws.Range("L54").Value = 18999
$ws.Range("N54").Value = -20039

$ws.Range("H136").Value = 4607.7354
$ws.Range("I136").Value = 2103
$ws.Range("J136").Value = 7597.2583
$ws.Range("K136").Value = 6309
$ws.Range("L136").Value = 22791.7749
$ws.Range("M136").Value = -3759
$ws.Range("N136").Value = -27891.7749
